$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 978.1429000000001
$ws.Range("J17").Value = 978.1429000000001
$ws.Range("L17").Value = 2934.4287
$ws.Range("N17").Value = -3270.4287
$ws.Range("H28").Value = 1610.5
$ws.Range("I28").Value = 2082.1428
$ws.Range("J28").Value = 1138.8572
$ws.Range("K28").Value = 2082.1428
$ws.Range("L28").Value = 1138.8572
$ws.Range("M28").Value = -1597.1428
$ws.Range("N28").Value = -2108.8572
$ws.Range("H62").Value = 1837.8572
$ws.Range("I62").Value = 1444.0952
$ws.Range("J62").Value = 2231.6191
$ws.Range("K62").Value = 1444.0952
$ws.Range("L62").Value = 2231.6191
$ws.Range("M62").Value = -820.0952
$ws.Range("N62").Value = -3479.6191
$ws.Range("H65").Value = 1837.8572
$ws.Range("I65").Value = 1444.0952
$ws.Range("J65").Value = 2231.6191
$ws.Range("K65").Value = 7220.476
$ws.Range("L65").Value = 11158.0955
$ws.Range("M65").Value = -4100.476
$ws.Range("N65").Value = -17398.0955
$ws.Range("H98").Value = 824.7241
$ws.Range("I98").Value = 837.43475
$ws.Range("J98").Value = 776
$ws.Range("K98").Value = 837.43475
$ws.Range("L98").Value = 776
$ws.Range("M98").Value = 660.56525
$ws.Range("N98").Value = -3772
$ws.Range("H112").Value = 2181.3157
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 2343.8235
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 7031.470499999999
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -9247.470499999999
$ws.Range("H122").Value = 824.7241
$ws.Range("I122").Value = 837.43475
$ws.Range("J122").Value = 776
$ws.Range("K122").Value = 2512.30425
$ws.Range("L122").Value = 2328
$ws.Range("M122").Value = -62.30425000000014
$ws.Range("N122").Value = -7228
$ws.Range("H132").Value = 767729.9
$ws.Range("I132").Value = 2049.64
$ws.Range("K132").Value = 6148.92
$ws.Range("M132").Value = -3618.92
$ws.Range("H135").Value = 20783.434
$ws.Range("I135").Value = 25672.244
$ws.Range("K135").Value = 231050.196
$ws.Range("M135").Value = -228515.196
$ws.Range("H137").Value = 1668030.1
$ws.Range("I137").Value = 2223167
$ws.Range("J137").Value = 2619.1333
$ws.Range("K137").Value = 6669501
$ws.Range("L137").Value = 7857.3999
$ws.Range("M137").Value = -6666951
$ws.Range("N137").Value = -12957.3999
$ws.Range("H138").Value = 2224201.2
$ws.Range("I138").Value = 1450.2954
$ws.Range("J138").Value = 5379073.5
$ws.Range("K138").Value = 4350.8862
$ws.Range("L138").Value = 16137220.5
$ws.Range("M138").Value = 789.1138000000001
$ws.Range("N138").Value = -16147500.5
$ws.Range("H141").Value = 2600.5715
$ws.Range("I141").Value = 1363.5682
$ws.Range("J141").Value = 7136.25
$ws.Range("K141").Value = 4090.7046
$ws.Range("L141").Value = 21408.75
$ws.Range("M141").Value = 1089.2954
$ws.Range("N141").Value = -31768.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2945.34
$ws.Range("I32").Value = 2551.8271
$ws.Range("K32").Value = 2551.8271
$ws.Range("M32").Value = -2264.8271
$ws.Range("H36").Value = 8271
$ws.Range("I36").Value = 8081.5
$ws.Range("J36").Value = 9029
$ws.Range("K36").Value = 8081.5
$ws.Range("L36").Value = 9029
$ws.Range("M36").Value = -7735.5
$ws.Range("N36").Value = -9721
$ws.Range("H61").Value = 31314362
$ws.Range("I61").Value = 43523204
$ws.Range("J61").Value = 113992.445
$ws.Range("K61").Value = 43523204
$ws.Range("L61").Value = 113992.445
$ws.Range("M61").Value = -43522992
$ws.Range("N61").Value = -114416.445
$ws.Range("H136").Value = 31314362
$ws.Range("I136").Value = 43523204
$ws.Range("J136").Value = 113992.445
$ws.Range("K136").Value = 130569612
$ws.Range("L136").Value = 341977.335
$ws.Range("M136").Value = -130567062
$ws.Range("N136").Value = -347077.335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 16668834
$ws.Range("I105").Value = 25002026
$ws.Range("J105").Value = 2450
$ws.Range("K105").Value = 25002026
$ws.Range("L105").Value = 2450
$ws.Range("M105").Value = -25000279
$ws.Range("N105").Value = -5944

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2026
$ws.Range("I31").Value = 1130.1333
$ws.Range("K31").Value = 1130.1333
$ws.Range("M31").Value = -835.1333
$ws.Range("H34").Value = 2026
$ws.Range("I34").Value = 1130.1333
$ws.Range("K34").Value = 1130.1333
$ws.Range("M34").Value = -928.1333
$ws.Range("H107").Value = 413.58823
$ws.Range("I107").Value = 376.125
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 376.125
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = 1543.875
$ws.Range("N107").Value = -4853
$ws.Range("H132").Value = 14881.243
$ws.Range("I132").Value = 1160.6034
$ws.Range("J132").Value = 64618.562
$ws.Range("K132").Value = 3481.8102
$ws.Range("L132").Value = 193855.686
$ws.Range("M132").Value = -951.8101999999999
$ws.Range("N132").Value = -198915.686
$ws.Range("H134").Value = 17864.215
$ws.Range("I134").Value = 1254.9111
$ws.Range("J134").Value = 55235.15
$ws.Range("K134").Value = 3764.7333
$ws.Range("L134").Value = 165705.45
$ws.Range("M134").Value = -1229.7333
$ws.Range("N134").Value = -170775.45

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7145273
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 7145273
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 21435819
$ws.Range("N4").Value = -21436043
$ws.Range("H5").Value = 568.17645
$ws.Range("I5").Value = 496.93332
$ws.Range("J5").Value = 1102.5
$ws.Range("K5").Value = 1490.79996
$ws.Range("L5").Value = 3307.5
$ws.Range("M5").Value = -1378.79996
$ws.Range("N5").Value = -3531.5
$ws.Range("H12").Value = 67.382355
$ws.Range("I12").Value = 79.666664
$ws.Range("J12").Value = 57.68421
$ws.Range("K12").Value = 238.999992
$ws.Range("L12").Value = 173.05263
$ws.Range("M12").Value = -65.99999199999999
$ws.Range("N12").Value = -519.05263
$ws.Range("H14").Value = 1006.8182
$ws.Range("I14").Value = 1006.8182
$ws.Range("K14").Value = 3020.4546
$ws.Range("M14").Value = -2847.4546
$ws.Range("H56").Value = 169914.6
$ws.Range("I56").Value = 169914.6
$ws.Range("K56").Value = 169914.6
$ws.Range("M56").Value = -169384.6
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H107").Value = 740.42426
$ws.Range("J107").Value = 704.3158
$ws.Range("L107").Value = 2112.9474
$ws.Range("N107").Value = -5952.9474
$ws.Range("H121").Value = 42886824
$ws.Range("I121").Value = 1388.3334
$ws.Range("J121").Value = 49841216
$ws.Range("K121").Value = 4165.0002
$ws.Range("L121").Value = 149523648
$ws.Range("M121").Value = -2855.0002
$ws.Range("N121").Value = -149526268
$ws.Range("H122").Value = 619.9231
$ws.Range("I122").Value = 282.4737
$ws.Range("J122").Value = 940.5
$ws.Range("K122").Value = 2542.2633
$ws.Range("L122").Value = 8464.5
$ws.Range("M122").Value = -92.26330000000007
$ws.Range("N122").Value = -13364.5
$ws.Range("H131").Value = 1039.2877
$ws.Range("J131").Value = 1115.9692
$ws.Range("L131").Value = 3347.9076
$ws.Range("N131").Value = -13427.9076
$ws.Range("H135").Value = 568.17645
$ws.Range("I135").Value = 496.93332
$ws.Range("J135").Value = 1102.5
$ws.Range("K135").Value = 4472.39988
$ws.Range("L135").Value = 9922.5
$ws.Range("M135").Value = -1937.39988
$ws.Range("N135").Value = -14992.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2204.6
$ws.Range("I122").Value = 2103.5
$ws.Range("J122").Value = 2272
$ws.Range("K122").Value = 6310.5
$ws.Range("L122").Value = 6816
$ws.Range("M122").Value = -3860.5
$ws.Range("N122").Value = -11716
$ws.Range("H126").Value = 2429.8125
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 2761.5454
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 8284.636200000001
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -13224.6362
$ws.Range("H132").Value = 54398.656
$ws.Range("I132").Value = 33778.613
$ws.Range("J132").Value = 145716
$ws.Range("K132").Value = 101335.839
$ws.Range("L132").Value = 437148
$ws.Range("M132").Value = -98805.83899999999
$ws.Range("N132").Value = -442208

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1075.1666
$ws.Range("I93").Value = 1075.1666
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1075.1666
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 172.8334
$ws.Range("H132").Value = 24651.637
$ws.Range("I132").Value = 1631.2858
$ws.Range("J132").Value = 114175.22
$ws.Range("K132").Value = 4893.857400000001
$ws.Range("L132").Value = 342525.66
$ws.Range("M132").Value = -2363.857400000001
$ws.Range("N132").Value = -347585.66
$ws.Range("H136").Value = 36810.945
$ws.Range("I136").Value = 26438.666
$ws.Range("J136").Value = 60606.176
$ws.Range("K136").Value = 79315.99800000001
$ws.Range("L136").Value = 181818.528
$ws.Range("M136").Value = -76765.99800000001
$ws.Range("N136").Value = -186918.528

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1895.6842
$ws.Range("I122").Value = 1321.6818
$ws.Range("J122").Value = 2684.9375
$ws.Range("K122").Value = 3965.0454
$ws.Range("L122").Value = 8054.8125
$ws.Range("M122").Value = -1515.0454
$ws.Range("N122").Value = -12954.8125
$ws.Range("H132").Value = 40034.5
$ws.Range("I132").Value = 28339.25
$ws.Range("J132").Value = 73449.5
$ws.Range("K132").Value = 85017.75
$ws.Range("L132").Value = 220348.5
$ws.Range("M132").Value = -82487.75
$ws.Range("H136").Value = 34713.785
$ws.Range("I136").Value = 31410.94
$ws.Range("J136").Value = 38606.43
$ws.Range("K136").Value = 94232.81999999999
$ws.Range("L136").Value = 115819.29
$ws.Range("M136").Value = -91682.81999999999
$ws.Range("N136").Value = -120919.29
